# Update CBC Blood Automated ValueSet metadata sheet:
#  - bump Version, Status, Date
#  - replace Contact with org contact + add a new person contact row
#  - add a new Jurisdiction row
#  - everything below shifts down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- simple field updates -------------------------------------------------
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-23T10:17:11-05:00"

# --- make room for the new "Contact" (person) row and the new
#     "Jurisdiction" row: shift rows 11..15 down to 12..16 -----------------
for ($r = 15; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Range("A$destRow").Value = $ws.Range("A$r").Value()
    $ws.Range("B$destRow").Value = $ws.Range("B$r").Value()
}

# Row 16 did not exist before, so it has no formatting yet - copy it from
# the row above (which carries the correct style) before writing new text.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Contact (organisation) now includes the CIBMTR url -------------------
$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- new Contact (person) row ----------------------------------------------
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- new Jurisdiction row ---------------------------------------------------
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Rows 13-16 already hold the right text after the shift above
#  13: Description / CBC panel - Blood by Automated count (58410-2)
#  14: Purpose / (blank)
#  15: Copyright / (blank)
#  16: Immutable / BooleanType[null]
